# Applies the "javaVScpp of CPT/+" update: adds query-time / memory columns
# for Java vs cpp CPT/CPT+ comparisons, extends the formulas area (J/K) as
# shared formulas down through row 28, and updates the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New column widths for the additional columns T:AB (20-28) ---
$ws.Columns.Item(20).ColumnWidth = 29.830729166666668
$ws.Columns.Item(21).ColumnWidth = 27.330729166666668
$ws.Columns.Item(22).ColumnWidth = 21.330729166666668
$ws.Columns.Item(23).ColumnWidth = 24.830729166666668
$ws.Columns.Item(24).ColumnWidth = 23.998697916666668
$ws.Columns.Item(25).ColumnWidth = 25.330729166666668
$ws.Columns.Item(26).ColumnWidth = 28.830729166666668
$ws.Columns.Item(27).ColumnWidth = 27.830729166666668
$ws.Columns.Item(28).ColumnWidth = 29.830729166666668

# --- 2. New header cells in row 1 (U1:Y1, then Z1:AB1) and the footnote
#        marker in T19. The order mirrors how the strings were originally
#        typed in, so that the shared-string table is populated the same
#        way (U..Y first, then the footnote at T19, then Z..AB).
$ws.Range("U1").Value = "Java CPT query time (ms)"
$ws.Range("V1").Value = "Java CPT+ query time (ms)"
$ws.Range("W1").Value = "cpp CPT query time (ms)"
$ws.Range("X1").Value = "cpp CPT+ query time (ms)"
$ws.Range("Y1").Value = "Java CPT memory (MB)*"
$ws.Range("T19").Value = "*As reported with calculations"
$ws.Range("Z1").Value = "Java CPT+ memory (MB)*"
$ws.Range("AA1").Value = "cpp CPT memory (MB)*"
$ws.Range("AB1").Value = "cpp CPT+ memory (MB)*"

# --- 3. Extend the J (CPT+ time ratio) and K (SD_CPT/CPT+ time) formulas
#        down as shared formulas across rows 19:28, matching the existing
#        H/I/M/N shared-formula pattern already used in that block.
$ws.Range("J19:J28").Formula = "=(M2+L2+I2+G2)/F2"
$ws.Range("K19:K28").Formula = "=I19/H19"

# --- 4. Update the sheet view: scroll so column S is the left-most
#        visible column and select Z23.
$excel.ActiveWindow.ScrollColumn = 19
$ws.Range("Z23").Select()
